# "add in lifestage columns, #855"
#
# The "life Stage" header (cell AC1 on the Template sheet) is renamed to
# "Lifestage". Re-assigning the cell's value makes Excel drop the old
# shared-string entry and append the new text to the shared-string table,
# which is exactly what the canonical diff shows (the "life Stage" <si>
# is removed and a new "Lifestage" <si> is appended at the end; every
# other header keeps its own text, just re-pointed at shifted shared-
# string indices as a consequence).
#
# The cell that was edited ends up selected/active, same as it would after
# a person typed the new header text in on screen and left the cursor
# there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AC1").Value = "Lifestage"

$ws.Range("AC1").Select()
